$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 180.33333
$ws.Range("I15").Value = 180.33333
$ws.Range("K15").Value = 540.99999
$ws.Range("M15").Value = -371.99999

$ws.Range("H33").Value = 532
$ws.Range("I33").Value = 286.125
$ws.Range("K33").Value = 286.125
$ws.Range("M33").Value = -57.125

$ws.Range("H53").Value = 644.7646999999999
$ws.Range("I53").Value = 574.53845
$ws.Range("K53").Value = 574.53845
$ws.Range("M53").Value = 62.46154999999999

$ws.Range("H76").Value = 3799.625
$ws.Range("I76").Value = 3699.5715
$ws.Range("K76").Value = 3699.5715
$ws.Range("M76").Value = -3384.5715

$ws.Range("H79").Value = 3799.625
$ws.Range("I79").Value = 3699.5715
$ws.Range("K79").Value = 3699.5715
$ws.Range("M79").Value = -2607.5715

$ws.Range("H80").Value = 3000
$ws.Range("I80").Value = 1500
$ws.Range("J80").Value = 4200
$ws.Range("K80").Value = 4500
$ws.Range("L80").Value = 12600
$ws.Range("M80").Value = -3502
$ws.Range("N80").Value = -14596

$ws.Range("H83").Value = 3000
$ws.Range("I83").Value = 1500
$ws.Range("J83").Value = 4200
$ws.Range("K83").Value = 13500
$ws.Range("L83").Value = 37800
$ws.Range("M83").Value = -8508
$ws.Range("N83").Value = -47784

$ws.Range("H92").Value = 1246.25
$ws.Range("I92").Value = 1246.25
$ws.Range("K92").Value = 1246.25
$ws.Range("M92").Value = 1.75

$ws.Range("H129").Value = 5348.5
$ws.Range("I129").Value = 5348.5
$ws.Range("K129").Value = 16045.5
$ws.Range("M129").Value = -11045.5

$ws.Range("H131").Value = 1281.6666
$ws.Range("I131").Value = 1281.6666
$ws.Range("K131").Value = 3844.9998
$ws.Range("M131").Value = 1195.0002

$ws.Range("H137").Value = 1146.4286
$ws.Range("I137").Value = 1005.6667
$ws.Range("J137").Value = 1991
$ws.Range("K137").Value = 3017.0001
$ws.Range("L137").Value = 5973
$ws.Range("M137").Value = -467.0001000000002
$ws.Range("N137").Value = -11073

$ws.Range("H138").Value = 4633
$ws.Range("I138").Value = 4633
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 13899
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -8759
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1482.909
$ws.Range("I2").Value = 551.6667
$ws.Range("K2").Value = 551.6667
$ws.Range("M2").Value = -438.6667

$ws.Range("H6").Value = 3000
$ws.Range("I6").Value = 3000
$ws.Range("J6").Value = 3000
$ws.Range("K6").Value = 3000
$ws.Range("L6").Value = 3000
$ws.Range("M6").Value = -2827
$ws.Range("N6").Value = -3346

$ws.Range("H32").Value = 4799.5713
$ws.Range("I32").Value = 3603.1765
$ws.Range("J32").Value = 9884.25
$ws.Range("K32").Value = 3603.1765
$ws.Range("L32").Value = 9884.25
$ws.Range("M32").Value = -3316.1765
$ws.Range("N32").Value = -10458.25

$ws.Range("H45").Value = 3500.389
$ws.Range("I45").Value = 1253
$ws.Range("J45").Value = 4142.5
$ws.Range("K45").Value = 1253
$ws.Range("L45").Value = 4142.5
$ws.Range("M45").Value = -876
$ws.Range("N45").Value = -4896.5

$ws.Range("H102").Value = 4066.6667
$ws.Range("I102").Value = 4066.6667
$ws.Range("K102").Value = 4066.6667
$ws.Range("M102").Value = -2444.6667

$ws.Range("H116").Value = 1482.909
$ws.Range("I116").Value = 551.6667
$ws.Range("K116").Value = 551.6667
$ws.Range("M116").Value = 1742.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1482.909
$ws.Range("I3").Value = 551.6667
$ws.Range("K3").Value = 551.6667
$ws.Range("M3").Value = -437.6667

$ws.Range("H33").Value = 1699
$ws.Range("I33").Value = 1699
$ws.Range("K33").Value = 1699
$ws.Range("M33").Value = -1363

$ws.Range("H39").Value = 55000
$ws.Range("J39").Value = 55000
$ws.Range("L39").Value = 55000
$ws.Range("N39").Value = -55778

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()

$ws.Range("H99").Value = 3394.524
$ws.Range("I99").Value = 3599.1667
$ws.Range("K99").Value = 3599.1667
$ws.Range("M99").Value = -2101.1667

$ws.Range("H107").Value = 305.5
$ws.Range("I107").Value = 305.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 305.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1614.5
$ws.Range("N107").ClearContents()

$ws.Range("H126").Value = 3394.524
$ws.Range("I126").Value = 3599.1667
$ws.Range("K126").Value = 10797.5001
$ws.Range("M126").Value = -8327.500100000001

$ws.Range("H132").Value = 4253.7144
$ws.Range("I132").Value = 4631.1665
$ws.Range("K132").Value = 13893.4995
$ws.Range("M132").Value = -11363.4995

$ws.Range("H134").Value = 2489.5715
$ws.Range("I134").Value = 2568.8333
$ws.Range("J134").Value = 2014
$ws.Range("K134").Value = 7706.499899999999
$ws.Range("L134").Value = 6042
$ws.Range("M134").Value = -5171.499899999999
$ws.Range("N134").Value = -11112

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 766.7143
$ws.Range("J12").Value = 766.7143
$ws.Range("L12").Value = 2300.1429
$ws.Range("N12").Value = -2646.1429

$ws.Range("H18").Value = 1842.7142
$ws.Range("I18").Value = 1449.5
$ws.Range("J18").Value = 2000
$ws.Range("K18").Value = 4348.5
$ws.Range("L18").Value = 6000
$ws.Range("M18").Value = -4179.5
$ws.Range("N18").Value = -6338

$ws.Range("H70").Value = 1970.6666
$ws.Range("I70").Value = 1970.6666
$ws.Range("K70").Value = 5911.9998
$ws.Range("M70").Value = -5596.9998

$ws.Range("H73").Value = 1970.6666
$ws.Range("I73").Value = 1970.6666
$ws.Range("K73").Value = 5911.9998
$ws.Range("M73").Value = -4819.9998

$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("M98").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 5411
$ws.Range("I107").Value = 472.85715
$ws.Range("J107").Value = 16933.334
$ws.Range("K107").Value = 472.85715
$ws.Range("L107").Value = 16933.334
$ws.Range("M107").Value = 1447.14285
$ws.Range("N107").Value = -20773.334

$ws.Range("H132").Value = 2764
$ws.Range("I132").Value = 2764
$ws.Range("K132").Value = 8292
$ws.Range("M132").Value = -5762

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H53").Value = 46499.5
$ws.Range("I53").Value = 46499.5
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 46499.5
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -45981.5
$ws.Range("N53").ClearContents()

$ws.Range("H61").Value = 706.5
$ws.Range("I61").Value = 667.8
$ws.Range("K61").Value = 667.8
$ws.Range("M61").Value = -465.8

$ws.Range("H93").Value = 612.55554
$ws.Range("I93").Value = 589.125
$ws.Range("J93").Value = 800
$ws.Range("K93").Value = 589.125
$ws.Range("L93").Value = 800
$ws.Range("M93").Value = 658.875
$ws.Range("N93").Value = -3296

$ws.Range("H113").Value = 706.5
$ws.Range("I113").Value = 667.8
$ws.Range("K113").Value = 667.8
$ws.Range("M113").Value = 1502.2

$ws.Range("H122").Value = 4751.0835
$ws.Range("I122").Value = 3891.6667
$ws.Range("K122").Value = 11675.0001
$ws.Range("M122").Value = -9225.000100000001

$ws.Range("H132").Value = 4100.5
$ws.Range("I132").Value = 3299.8333
$ws.Range("J132").Value = 6502.5
$ws.Range("K132").Value = 9899.499899999999
$ws.Range("L132").Value = 19507.5
$ws.Range("M132").Value = -7369.499899999999
$ws.Range("N132").Value = -24567.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 5000000
$ws.Range("J18").Value = 5000000
$ws.Range("L18").Value = 5000000
$ws.Range("N18").Value = -5000346

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()
